$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "x" mark that was in G22 to a "1gg" (1 day) estimate in J22/K22,
# reflecting the newly added hard-difficulty sudoku task.
$ws.Range("G22").ClearContents() | Out-Null
$ws.Range("J22").Value = "1gg"
$ws.Range("K22").Value = "1gg"

# Update the active cell selection to D25
$ws.Range("D25").Select() | Out-Null
